$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6 updates (ano = 2025):
# total_customers 416 -> 418
$ws.Range("C6").Value = 418
# new_customers 108 -> 110
$ws.Range("E6").Value = 110
# new_rate 25.96153846153846 -> 26.31578947368421
$ws.Range("G6").Value = 26.31578947368421
# returning_rate 74.03846153846155 -> 73.68421052631578
$ws.Range("H6").Value = 73.68421052631578
